$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in C3: "int jloop_ub=0" -> "int j_loop_ub=0"
$ws.Range("C3").Value = "int j_loop_ub=0"

# Update the selection to C7 (as recorded in the saved view state)
$ws.Range("C7").Select()
